# Elimina EC anteriores y se agregan nuevos, se modifica base de datos
#
# The "Periodo Mora" detail table (rows 16-56) mixed periods for two
# workers (CARLOS HUMBERTO VILLA MERCADO and ESTEBAN DARIO BARBOZA PRIMERA).
# The old data had Carlos's 4 periods (2104-2107) interleaved with
# Esteban's block, and Esteban's block duplicated those same 4 periods.
# The new data de-duplicates: Carlos's 4 periods now form a clean block at
# the top (rows 16-19, newest period first), followed by Esteban's full
# period history (rows 20-56, newest period first, 2104-2107 removed since
# those now exclusively belong to Carlos).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @{Row=16; C="1143362651"; D="CARLOS HUMBERTO VILLA MERCADO"; E="2107"; F=580000;  G=17127677},
    @{Row=17; C="1143362651"; D="CARLOS HUMBERTO VILLA MERCADO"; E="2106"; F=580000;  G=17127677},
    @{Row=18; C="1143362651"; D="CARLOS HUMBERTO VILLA MERCADO"; E="2105"; F=580000;  G=17127677},
    @{Row=19; C="1143362651"; D="CARLOS HUMBERTO VILLA MERCADO"; E="2104"; F=580000;  G=17127677},
    @{Row=20; C="1143380298"; D="ESTEBAN DARIO BARBOZA PRIMERA"; E="2312"; F=257133;  G=14500000},
    @{Row=21; C="1143380298"; D="ESTEBAN DARIO BARBOZA PRIMERA"; E="2311"; F=406000;  G=14500000},
    @{Row=22; C="1143380298"; D="ESTEBAN DARIO BARBOZA PRIMERA"; E="2310"; F=406000;  G=14500000},
    @{Row=23; C="1143380298"; D="ESTEBAN DARIO BARBOZA PRIMERA"; E="2309"; F=406000;  G=14500000},
    @{Row=24; C="1143380298"; D="ESTEBAN DARIO BARBOZA PRIMERA"; E="2308"; F=406000;  G=14500000},
    @{Row=25; C="1143380298"; D="ESTEBAN DARIO BARBOZA PRIMERA"; E="2307"; F=406000;  G=14500000},
    @{Row=26; C="1143380298"; D="ESTEBAN DARIO BARBOZA PRIMERA"; E="2306"; F=406000;  G=14500000},
    @{Row=27; C="1143380298"; D="ESTEBAN DARIO BARBOZA PRIMERA"; E="2305"; F=406000;  G=14500000},
    @{Row=28; C="1143380298"; D="ESTEBAN DARIO BARBOZA PRIMERA"; E="2304"; F=406000;  G=14500000},
    @{Row=29; C="1143380298"; D="ESTEBAN DARIO BARBOZA PRIMERA"; E="2303"; F=406000;  G=14500000},
    @{Row=30; C="1143380298"; D="ESTEBAN DARIO BARBOZA PRIMERA"; E="2302"; F=406000;  G=14500000},
    @{Row=31; C="1143380298"; D="ESTEBAN DARIO BARBOZA PRIMERA"; E="2301"; F=406000;  G=14500000},
    @{Row=32; C="1143380298"; D="ESTEBAN DARIO BARBOZA PRIMERA"; E="2212"; F=406000;  G=14500000},
    @{Row=33; C="1143380298"; D="ESTEBAN DARIO BARBOZA PRIMERA"; E="2211"; F=406000;  G=14500000},
    @{Row=34; C="1143380298"; D="ESTEBAN DARIO BARBOZA PRIMERA"; E="2210"; F=406000;  G=14500000},
    @{Row=35; C="1143380298"; D="ESTEBAN DARIO BARBOZA PRIMERA"; E="2209"; F=406000;  G=14500000},
    @{Row=36; C="1143380298"; D="ESTEBAN DARIO BARBOZA PRIMERA"; E="2208"; F=406000;  G=14500000},
    @{Row=37; C="1143380298"; D="ESTEBAN DARIO BARBOZA PRIMERA"; E="2207"; F=406000;  G=14500000},
    @{Row=38; C="1143380298"; D="ESTEBAN DARIO BARBOZA PRIMERA"; E="2206"; F=406000;  G=14500000},
    @{Row=39; C="1143380298"; D="ESTEBAN DARIO BARBOZA PRIMERA"; E="2205"; F=406000;  G=14500000},
    @{Row=40; C="1143380298"; D="ESTEBAN DARIO BARBOZA PRIMERA"; E="2204"; F=406000;  G=14500000},
    @{Row=41; C="1143380298"; D="ESTEBAN DARIO BARBOZA PRIMERA"; E="2203"; F=406000;  G=14500000},
    @{Row=42; C="1143380298"; D="ESTEBAN DARIO BARBOZA PRIMERA"; E="2202"; F=406000;  G=14500000},
    @{Row=43; C="1143380298"; D="ESTEBAN DARIO BARBOZA PRIMERA"; E="2201"; F=406000;  G=14500000},
    @{Row=44; C="1143380298"; D="ESTEBAN DARIO BARBOZA PRIMERA"; E="2112"; F=406000;  G=14500000},
    @{Row=45; C="1143380298"; D="ESTEBAN DARIO BARBOZA PRIMERA"; E="2111"; F=406000;  G=14500000},
    @{Row=46; C="1143380298"; D="ESTEBAN DARIO BARBOZA PRIMERA"; E="2110"; F=406000;  G=14500000},
    @{Row=47; C="1143380298"; D="ESTEBAN DARIO BARBOZA PRIMERA"; E="2109"; F=406000;  G=14500000},
    @{Row=48; C="1143380298"; D="ESTEBAN DARIO BARBOZA PRIMERA"; E="2108"; F=406000;  G=14500000},
    @{Row=49; C="1143380298"; D="ESTEBAN DARIO BARBOZA PRIMERA"; E="2107"; F=406000;  G=14500000},
    @{Row=50; C="1143380298"; D="ESTEBAN DARIO BARBOZA PRIMERA"; E="2106"; F=406000;  G=14500000},
    @{Row=51; C="1143380298"; D="ESTEBAN DARIO BARBOZA PRIMERA"; E="2105"; F=406000;  G=14500000},
    @{Row=52; C="1143380298"; D="ESTEBAN DARIO BARBOZA PRIMERA"; E="2104"; F=406000;  G=14500000},
    @{Row=53; C="1143380298"; D="ESTEBAN DARIO BARBOZA PRIMERA"; E="2103"; F=406000;  G=14500000},
    @{Row=54; C="1143380298"; D="ESTEBAN DARIO BARBOZA PRIMERA"; E="2102"; F=406000;  G=14500000},
    @{Row=55; C="1143380298"; D="ESTEBAN DARIO BARBOZA PRIMERA"; E="2101"; F=406000;  G=14500000},
    @{Row=56; C="1143380298"; D="ESTEBAN DARIO BARBOZA PRIMERA"; E="2012"; F=230066;  G=14500000}
)

foreach ($r in $rows) {
    $ws.Range("C$($r.Row)").Value = $r.C
    $ws.Range("D$($r.Row)").Value = $r.D
    $ws.Range("E$($r.Row)").Value = $r.E
    $ws.Range("F$($r.Row)").Value = $r.F
    $ws.Range("G$($r.Row)").Value = $r.G
}
